$d = $word.ActiveDocument

function Set-ParagraphXml($paragraphIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paragraphIndex)
    $range = $d.Range($p.Range.Start, $p.Range.End)
    $full = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($full)
}

Set-ParagraphXml 37 '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">If you ever need to access an old file version for the sake of functionality, feel free to </w:t></w:r><w:r><w:t xml:space="preserve">click on the file in question, and </w:t></w:r><w:r><w:t xml:space="preserve">look at </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>History</w:t></w:r><w:r><w:t xml:space="preserve"> for any one file</w:t></w:r><w:r><w:t xml:space="preserve">. Clicking </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>View at this point in history</w:t></w:r><w:r><w:t xml:space="preserve"> allows you to download the file at any point it was uploaded here,</w:t></w:r><w:r><w:t xml:space="preserve"> and you can locally restore original functionality to the </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>S:\</w:t></w:r><w:r><w:t xml:space="preserve"> drive folder</w:t></w:r><w:r><w:t xml:space="preserve"> in the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>off chance</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> functionality is impeded</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>'
Set-ParagraphXml 36 '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">This will commit the changes to the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>repository, and</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> create a </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>new version</w:t></w:r><w:r><w:t xml:space="preserve"> of the specified file.</w:t></w:r></w:p>'
Set-ParagraphXml 35 '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Choose the </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>file</w:t></w:r><w:r><w:t xml:space="preserve"> you just changed (GitHub will automatically recognize any </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>file changes</w:t></w:r><w:r><w:t xml:space="preserve">, but does not allow </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>folder uploads</w:t></w:r><w:r><w:t xml:space="preserve">), </w:t></w:r><w:r><w:t xml:space="preserve">enter a </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">description of the </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>change</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>and</w:t></w:r><w:r><w:t xml:space="preserve"> press the green button near the bottom named </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Commit changes</w:t></w:r></w:p>'
Set-ParagraphXml 22 '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="red"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="red"/></w:rPr><w:t>VFD_Template</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="red"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:highlight w:val="red"/></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:highlight w:val="red"/></w:rPr><w:t xml:space="preserve"> Backup</w:t></w:r></w:p>'
Set-ParagraphXml 20 '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="red"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="red"/></w:rPr><w:t>Emailing_Saving</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="red"/></w:rPr><w:t xml:space="preserve"> Macros</w:t></w:r></w:p>'
Set-ParagraphXml 19 '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="red"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="red"/></w:rPr><w:t>DriveSelect</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="red"/></w:rPr><w:t xml:space="preserve"> - No PDF Converter</w:t></w:r></w:p>'
Set-ParagraphXml 8 '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Motor_BOM_Tool</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Set-ParagraphXml 7 '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>M_Drives</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Set-ParagraphXml 6 '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>VFD_Template</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Set-ParagraphXml 5 '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>DriveSelect</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
